$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: replace "November" entry with "Someone Else"
$ws.Range("A2").Value = "Someone Else"
$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "5.9"
$ws.Range("B2").ClearFormats()
$ws.Range("C2").Value = "Col Spector"
$ws.Range("D2").Value = "Cameo Appearance of Neville Bardoliwalla as Freddie Mercury."

# Row 3: replace "Outlaw King" entry with "A Star Is Born"
$ws.Range("A3").Value = "A Star Is Born"
$ws.Range("B3").NumberFormat = "@"
$ws.Range("B3").Value = "7.8"
$ws.Range("B3").ClearFormats()
$ws.Range("C3").Value = "George Cukor"
$ws.Range("D3").Value = "A film star helps a young singer and actress find fame, even as age and alcoholism send his own career on a downward spiral."

# Row 4: delete entire row ("The Nun" entry removed)
$ws.Range("A4:D4").EntireRow.Delete()
